$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the FECHA (date) column from text strings to real date values
# for the rows that still stored them as text (rows 4-6), reusing the
# existing date number format/style already used by F2:F3.
$ws.Range("F4").Value = 44232   # 05/02/2021
$ws.Range("F5").Value = 44328   # 12/05/2021
$ws.Range("F6").Value = 44361   # 14/06/2021

$ws.Range("F2").Copy()
$ws.Range("F4:F6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add the new STATUS / COMMENTS columns
$ws.Range("I1").Value = "STATUS"
$ws.Range("I2").Value = "On time"
$ws.Range("J1").Value = "COMMENTS"
$ws.Range("J2").Value = "No comments yet"
$ws.Range("I3").Value = "On time"
$ws.Range("J3").Value = "No comments yet"
$ws.Range("I4").Value = "On time"
$ws.Range("J4").Value = "No comments yet"
$ws.Range("I5").Value = "On time"
$ws.Range("J5").Value = "No comments yet"
$ws.Range("I6").Value = "On time"
$ws.Range("J6").Value = "No comments yet"

# Give the new headers the same look as the existing ones
$ws.Range("A1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("J1").HorizontalAlignment = -4131

# Widen the new COMMENTS column a bit (closest width reachable through
# the engine's column-width quantization to the authored 18.5546875)
$ws.Range("J1").ColumnWidth = 17.666666666666668

# Match the selection left behind after the edit
$ws.Range("J7:K10").Select()
